$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: swap "Autor"/"Fecha" columns and rename "Imagen" -> "Imagenes" ---
$ws.Range("D1").Value = "Fecha"
$ws.Range("E1").Value = "Autor"
$ws.Range("F1").Value = "Imagenes"

# --- Row 2 ---
$ws.Range("D2").Value = "March 25, 2025"
$ws.Range("E2").Value = "Haik Aftandilian"
$ws.Range("F2").Value = "No se encontraron imagenes en este articulo"

# --- Row 3 ---
$ws.Range("D3").Value = "February 13, 2025"
$ws.Range("E3").Value = "James Graham"
$ws.Range("F3").ClearContents()

# --- Row 4 ---
$ws.Range("D4").Value = "December 4, 2024"
$ws.Range("E4").Value = "Mark Mayo"
$ws.Range("F4").Value = "No se encontraron imagenes en este articulo"

# --- Row 5 ---
$ws.Range("D5").Value = "October 16, 2024"
$ws.Range("E5").Value = "Stephen Hood"
$ws.Range("F5").ClearContents()

# --- Row 6 ---
$ws.Range("D6").Value = "August 8, 2024"
$ws.Range("E6").Value = "Marco Figueroa"
$ws.Range("F6").ClearContents()

# --- Row 7 ---
$ws.Range("D7").Value = "August 7, 2024"
$ws.Range("E7").Value = "James Graham"
$ws.Range("F7").Value = "No se encontraron imagenes en este articulo"

# --- Row 8 ---
$ws.Range("D8").Value = "June 27, 2024"
$ws.Range("E8").Value = "Christian Holler"
$ws.Range("F8").Value = "No se encontraron imagenes en este articulo"

# --- Row 9 ---
$ws.Range("D9").Value = "June 25, 2024"
$ws.Range("E9").Value = "Stephen Hood"
$ws.Range("F9").ClearContents()

# --- Row 10 ---
$ws.Range("D10").Value = "May 31, 2024"
$ws.Range("E10").Value = "Tarek Ziadé"
$ws.Range("F10").ClearContents()

# --- Row 11 ---
$ws.Range("D11").Value = "April 25, 2024"
$ws.Range("E11").Value = "Stephen Hood"
$ws.Range("F11").Value = "No se encontraron imagenes en este articulo"

# --- Row 12 ---
$ws.Range("D12").Value = "April 23, 2024"
$ws.Range("E12").Value = "Alex Franchuk"
$ws.Range("F12").ClearContents()
